$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Förändrad" date value (Excel serial 45807 == 2025-05-30), replacing the
# previous value of 45804 (2025-05-27) for every data row in column C.
$newDate = Get-Date -Year 2025 -Month 5 -Day 30 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
